{"js": "// Replace the date line and each multiplication-expression cell's text\n// with its updated value, per the authored diff. Every source string is\n// unique within the document, so a direct search+replace per pair is\n// unambiguous and keeps each run's formatting (font/size) untouched.\nconst replacements = [\n  [\"2025-01-23 Thursday\", \"2025-01-24 Friday\"],\n  [\"120\u00d75=\", \"420\u00d72=\"],\n  [\"345\u00d79=\", \"663\u00d79=\"],\n  [\"316\u00d73=\", \"189\u00d72=\"],\n  [\"219\u00d78=\", \"985\u00d77=\"],\n  [\"654\u00d77=\", \"711\u00d75=\"],\n  [\"564\u00d79=\", \"924\u00d73=\"],\n  [\"491\u00d79=\", \"707\u00d76=\"],\n  [\"207\u00d75=\", \"893\u00d79=\"],\n  [\"525\u00d73=\", \"483\u00d75=\"],\n  [\"115\u00d75=\", \"550\u00d73=\"],\n  [\"896\u00d75=\", \"248\u00d75=\"],\n  [\"839\u00d79=\", \"652\u00d79=\"],\n  [\"184\u00d74=\", \"690\u00d77=\"],\n  [\"442\u00d79=\", \"891\u00d74=\"],\n  [\"332\u00d74=\", \"660\u00d79=\"],\n  [\"249\u00d73=\", \"915\u00d72=\"],\n  [\"361\u00d75=\", \"698\u00d78=\"],\n  [\"692\u00d75=\", \"836\u00d74=\"],\n  [\"911\u00d78=\", \"771\u00d75=\"],\n  [\"801\u00d74=\", \"390\u00d73=\"],\n  [\"782\u00d72=\", \"456\u00d78=\"],\n  [\"165\u00d73=\", \"194\u00d74=\"],\n  [\"351\u00d79=\", \"435\u00d74=\"],\n  [\"264\u00d76=\", \"320\u00d73=\"],\n  [\"101\u00d76=\", \"991\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each multiplication-expression cell's text\n# with its updated value, per the authored diff. Every source string is\n# unique within the document, so Find/Replace per pair is unambiguous\n# and preserves each run's formatting (font/size).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-01-23 Thursday\", \"2025-01-24 Friday\"),\n    @(\"120\u00d75=\", \"420\u00d72=\"),\n    @(\"345\u00d79=\", \"663\u00d79=\"),\n    @(\"316\u00d73=\", \"189\u00d72=\"),\n    @(\"219\u00d78=\", \"985\u00d77=\"),\n    @(\"654\u00d77=\", \"711\u00d75=\"),\n    @(\"564\u00d79=\", \"924\u00d73=\"),\n    @(\"491\u00d79=\", \"707\u00d76=\"),\n    @(\"207\u00d75=\", \"893\u00d79=\"),\n    @(\"525\u00d73=\", \"483\u00d75=\"),\n    @(\"115\u00d75=\", \"550\u00d73=\"),\n    @(\"896\u00d75=\", \"248\u00d75=\"),\n    @(\"839\u00d79=\", \"652\u00d79=\"),\n    @(\"184\u00d74=\", \"690\u00d77=\"),\n    @(\"442\u00d79=\", \"891\u00d74=\"),\n    @(\"332\u00d74=\", \"660\u00d79=\"),\n    @(\"249\u00d73=\", \"915\u00d72=\"),\n    @(\"361\u00d75=\", \"698\u00d78=\"),\n    @(\"692\u00d75=\", \"836\u00d74=\"),\n    @(\"911\u00d78=\", \"771\u00d75=\"),\n    @(\"801\u00d74=\", \"390\u00d73=\"),\n    @(\"782\u00d72=\", \"456\u00d78=\"),\n    @(\"165\u00d73=\", \"194\u00d74=\"),\n    @(\"351\u00d79=\", \"435\u00d74=\"),\n    @(\"264\u00d76=\", \"320\u00d73=\"),\n    @(\"101\u00d76=\", \"991\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
